$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Parkinson" column header (I1), copying the header format
# (bold, centered, bordered) from the existing H1 header cell, then set
# the text.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("I1").Value = "Parkinson"

# Add the data values for the new "Parkinson" column (all zero, plain
# unstyled numeric cells like the rest of the data columns).
$ws.Range("I2").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("I4").Value = 0
